# Append a new daily row (2025-11-08) to the "Chart" sheet, mirroring the
# pattern of the existing rows (Date as text, Non-HTTPS count, HTTPS count).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newRow = 34

# Force the date-looking string to be stored as literal text (shared string)
# instead of being auto-converted to a date serial number. We briefly mark
# the cell as Text, assign the value, then clear the formatting again so the
# cell keeps using the default/general style like all the other rows.
$cell = $ws.Cells.Item($newRow, 1)
$cell.NumberFormat = "@"
$cell.Value = "2025-11-08"
$cell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 0.0
$ws.Cells.Item($newRow, 3).Value = 86.0
